$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "52.030.15"
$ws.Range("E2").Value = "  +1.16%  "

# Row 3
$ws.Range("D3").Value = "2.879.97"
$ws.Range("E3").Value = "  +3.57%  "

# Row 4
$ws.Range("E4").Value = "  +0.08%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "351.93"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.29%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "111.37"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +3.12%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.558"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +1.58%  "

# Row 8
$ws.Range("E8").Value = "  +0.07%  "

# Row 9
$ws.Range("E9").Value = "  +0.26%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "40.03"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +2.18%  "

# Row 11
$ws.Range("E11").Value = "  +3.52%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.135"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +0.32%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.03"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +0.69%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.80"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +0.75%  "

# Row 15
$ws.Range("D15").Value = "3.336.49"
$ws.Range("E15").Value = "  +3.76%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.993"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +7.38%  "

# Row 17
$ws.Range("D17").Value = "2.887.66"
$ws.Range("E17").Value = "  +3.41%  "

# Row 18
$ws.Range("D18").Value = "52.066.79"
$ws.Range("E18").Value = "  +1.24%  "

# Row 19
$ws.Range("B19").Value = "Uniswap"
$ws.Range("C19").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.73"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +0.34%  "

# Row 20
$ws.Range("B20").Value = "ImmutableX"
$ws.Range("C20").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "3.34"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +7.41%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.92"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +4.21%  "

# Row 22
$ws.Range("D22").Value = "0.0₃0981"
$ws.Range("E22").Value = "  +1.62%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "70.97"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +0.68%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "270.50"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +1.52%  "

# Row 25
$ws.Range("E25").Value = "  +0.43%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.30"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +1.75%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.999"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -0.11%  "

# Row 28
$ws.Range("E28").Value = "  +0.06%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "10.53"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +2.60%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "38.61"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +4.22%  "

# Row 31
$ws.Range("E31").Value = "  +0.77%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.42"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +3.27%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "53.24"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +2.78%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0941"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +11.26%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.93"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +4.42%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0457"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +3.29%  "

# Row 37
$ws.Range("E37").Value = "  -0.01%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.30"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +6.26%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.58"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +0.41%  "

# Row 40
$ws.Range("E40").Value = "  +3.57%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.63"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +5.83%  "

# Row 42
$ws.Range("E42").Value = "  +2.74%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "22.50"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +2.27%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "121.68"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +1.35%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.21"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +1.37%  "

# Row 46
$ws.Range("E46").Value = "  +6.90%  "

# Row 47
$ws.Range("D47").Value = "2.192.43"
$ws.Range("E47").Value = "  +2.91%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.51"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +8.05%  "

# Row 49
$ws.Range("E49").Value = "  +19.25%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.949"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +6.45%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "5.52"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +0.74%  "

Write-Output "Updated cryptos list"
